$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the "Posted At" column keeps plain text values (e.g. "2026-02-10")
# instead of Excel auto-converting them into date serials.
$ws.Range("F2:F23").NumberFormat = "@"

$ws.Cells.Item(2, 1).Value = "SR DATA SCIENTIST, SMART MFG & AI"
$ws.Cells.Item(2, 2).Value = "Micron Technology"
$ws.Cells.Item(2, 3).Value = "Boise, ID, US USA"
$ws.Cells.Item(2, 4).Value = 26.7
$ws.Cells.Item(2, 5).Value = "Data Scientist, RAG, TensorFlow, PyTorch, XGBoost, Keras, OpenCV, BigQuery, FastAPI, Docker"
$ws.Cells.Item(2, 6).Value = "2026-02-10"
$ws.Cells.Item(2, 7).Value = "https://www.indeed.com/viewjob?jk=12881722cc760412"

$ws.Cells.Item(3, 1).Value = "Senior Data Scientist"
$ws.Cells.Item(3, 2).Value = "Micron Technology"
$ws.Cells.Item(3, 3).Value = "Boise, ID, US USA"
$ws.Cells.Item(3, 4).Value = 25.6
$ws.Cells.Item(3, 5).Value = "Data Scientist, RAG, TensorFlow, PyTorch, XGBoost, Keras, BigQuery, FastAPI, Docker, Kubernetes"
$ws.Cells.Item(3, 6).Value = "2026-02-05"
$ws.Cells.Item(3, 7).Value = "https://www.indeed.com/viewjob?jk=9088b884cc79245f"

$ws.Cells.Item(4, 1).Value = "Senior Data Engineer, Business Intelligence"
$ws.Cells.Item(4, 2).Value = "Klaviyo"
$ws.Cells.Item(4, 3).Value = "Boston, MA, US USA"
$ws.Cells.Item(4, 4).Value = 15.6
$ws.Cells.Item(4, 5).Value = "RAG, S3, Data Lake, Docker, CI/CD, GitHub Actions, Terraform, Git, Snowflake, Python"
$ws.Cells.Item(4, 6).Value = "2026-02-27"
$ws.Cells.Item(4, 7).Value = "https://www.indeed.com/viewjob?jk=0f39474e07389a32"

$ws.Cells.Item(5, 1).Value = "MLOps Engineer - NO VISA SPONSORSHIP"
$ws.Cells.Item(5, 2).Value = "ProTalent Finders"
$ws.Cells.Item(5, 3).Value = "Fort Worth, TX, US USA"
$ws.Cells.Item(5, 4).Value = 15.6
$ws.Cells.Item(5, 5).Value = "AI Engineer, Data Scientist, Machine Learning Engineer, TensorFlow, PyTorch, Azure ML, Data Lake, MLflow, Kubernetes, CI/CD"
$ws.Cells.Item(5, 6).Value = "2026-02-27"
$ws.Cells.Item(5, 7).Value = "https://www.indeed.com/viewjob?jk=f7d5540e90f82629"

$ws.Cells.Item(6, 1).Value = "Senior Full Stack Developer"
$ws.Cells.Item(6, 2).Value = "Kalamata Capital Group, LLC."
$ws.Cells.Item(6, 3).Value = "New York, NY, US USA"
$ws.Cells.Item(6, 4).Value = 12.2
$ws.Cells.Item(6, 5).Value = "Docker, Kubernetes, CI/CD, GitHub Actions, Git, MongoDB, NoSQL, SQL, R, Java"
$ws.Cells.Item(6, 6).Value = "2026-02-27"
$ws.Cells.Item(6, 7).Value = "https://www.indeed.com/viewjob?jk=7e21a53fad2d2bec"

$ws.Cells.Item(7, 1).Value = "Software Engineer - Video"
$ws.Cells.Item(7, 2).Value = "Twilio"
$ws.Cells.Item(7, 3).Value = "Austin, TX, US USA"
$ws.Cells.Item(7, 4).Value = 11.1
$ws.Cells.Item(7, 5).Value = "RAG, Copilot, CI/CD, Git, Kafka, NoSQL, SQL, R, Java, Scala"
$ws.Cells.Item(7, 6).Value = "2026-02-27"
$ws.Cells.Item(7, 7).Value = "https://www.indeed.com/viewjob?jk=9a237a79dc97e081"

$ws.Cells.Item(8, 1).Value = "Technical SEO Engineer (Chicago, IL)"
$ws.Cells.Item(8, 2).Value = "nan"
$ws.Cells.Item(8, 3).Value = "Chicago, IL, US USA"
$ws.Cells.Item(8, 4).Value = 11.1
$ws.Cells.Item(8, 5).Value = "LangChain, RAG, CI/CD, Git, Python, SQL, R, Java, Scala, Optimization"
$ws.Cells.Item(8, 6).Value = "2026-02-27"
$ws.Cells.Item(8, 7).Value = "https://www.indeed.com/viewjob?jk=38710b2e8fbfdc85"

$ws.Cells.Item(9, 1).Value = "Technical SEO Engineer (Dallas, TX)"
$ws.Cells.Item(9, 2).Value = "nan"
$ws.Cells.Item(9, 3).Value = "Dallas, TX, US USA"
$ws.Cells.Item(9, 4).Value = 11.1
$ws.Cells.Item(9, 5).Value = "LangChain, RAG, CI/CD, Git, Python, SQL, R, Java, Scala, Optimization"
$ws.Cells.Item(9, 6).Value = "2026-02-27"
$ws.Cells.Item(9, 7).Value = "https://www.indeed.com/viewjob?jk=43bdd5164d993eaa"

$ws.Cells.Item(10, 1).Value = "Technical SEO Engineer (Indianapolis, IN)"
$ws.Cells.Item(10, 2).Value = "nan"
$ws.Cells.Item(10, 3).Value = "Indianapolis, IN, US USA"
$ws.Cells.Item(10, 4).Value = 11.1
$ws.Cells.Item(10, 5).Value = "LangChain, RAG, CI/CD, Git, Python, SQL, R, Java, Scala, Optimization"
$ws.Cells.Item(10, 6).Value = "2026-02-27"
$ws.Cells.Item(10, 7).Value = "https://www.indeed.com/viewjob?jk=96f0964c1340e086"

$ws.Cells.Item(11, 1).Value = "Technical SEO Engineer (Charlotte, NC)"
$ws.Cells.Item(11, 2).Value = "nan"
$ws.Cells.Item(11, 3).Value = "Charlotte, NC, US USA"
$ws.Cells.Item(11, 4).Value = 11.1
$ws.Cells.Item(11, 5).Value = "LangChain, RAG, CI/CD, Git, Python, SQL, R, Java, Scala, Optimization"
$ws.Cells.Item(11, 6).Value = "2026-02-27"
$ws.Cells.Item(11, 7).Value = "https://www.indeed.com/viewjob?jk=d65bb7b30a589154"

$ws.Cells.Item(12, 1).Value = "Technical SEO Engineer (Salk Lake City, UT)"
$ws.Cells.Item(12, 2).Value = "nan"
$ws.Cells.Item(12, 3).Value = "Salt Lake City, UT, US USA"
$ws.Cells.Item(12, 4).Value = 11.1
$ws.Cells.Item(12, 5).Value = "LangChain, RAG, CI/CD, Git, Python, SQL, R, Java, Scala, Optimization"
$ws.Cells.Item(12, 6).Value = "2026-02-27"
$ws.Cells.Item(12, 7).Value = "https://www.indeed.com/viewjob?jk=8022870adb0c65f9"

$ws.Cells.Item(13, 1).Value = "Technical SEO Engineer"
$ws.Cells.Item(13, 2).Value = "nan"
$ws.Cells.Item(13, 3).Value = "US USA"
$ws.Cells.Item(13, 4).Value = 11.1
$ws.Cells.Item(13, 5).Value = "LangChain, RAG, CI/CD, Git, Python, SQL, R, Java, Scala, Optimization"
$ws.Cells.Item(13, 6).Value = "2026-02-27"
$ws.Cells.Item(13, 7).Value = "https://www.indeed.com/viewjob?jk=4148a0a062852dba"

$ws.Cells.Item(14, 1).Value = "Software Engineer - Video"
$ws.Cells.Item(14, 2).Value = "Twilio"
$ws.Cells.Item(14, 3).Value = "US USA"
$ws.Cells.Item(14, 4).Value = 11.1
$ws.Cells.Item(14, 5).Value = "RAG, Copilot, CI/CD, Git, Kafka, NoSQL, SQL, R, Java, Scala"
$ws.Cells.Item(14, 6).Value = "2026-02-27"
$ws.Cells.Item(14, 7).Value = "https://www.indeed.com/viewjob?jk=e6902022dae079f0"

$ws.Cells.Item(15, 1).Value = "Python Engineer"
$ws.Cells.Item(15, 2).Value = "Realign"
$ws.Cells.Item(15, 3).Value = "Alpharetta, GA, US USA"
$ws.Cells.Item(15, 4).Value = 11.1
$ws.Cells.Item(15, 5).Value = "Data Scientist, TensorFlow, PyTorch, XGBoost, Docker, Kubernetes, CI/CD, Python, R, Optimization"
$ws.Cells.Item(15, 6).Value = "2026-02-27"
$ws.Cells.Item(15, 7).Value = "https://www.indeed.com/viewjob?jk=9031f7ae3ffe7ffe"

$ws.Cells.Item(16, 1).Value = "Data Engineer"
$ws.Cells.Item(16, 2).Value = "RogueSearch"
$ws.Cells.Item(16, 3).Value = "Remote, US USA"
$ws.Cells.Item(16, 4).Value = 11.1
$ws.Cells.Item(16, 5).Value = "Data Scientist, RAG, Data Lake, Kafka, Python, SQL, R, Java, Scala, Optimization"
$ws.Cells.Item(16, 6).Value = "2026-02-27"
$ws.Cells.Item(16, 7).Value = "https://www.indeed.com/viewjob?jk=0b30c7e8de15bdeb"

$ws.Cells.Item(17, 1).Value = "Data Scientist"
$ws.Cells.Item(17, 2).Value = "Bank of America"
$ws.Cells.Item(17, 3).Value = "Plano, TX, US USA"
$ws.Cells.Item(17, 4).Value = 11.1
$ws.Cells.Item(17, 5).Value = "Data Scientist, Generative AI, RAG, Git, Cassandra, NoSQL, Python, SQL, R, Java"
$ws.Cells.Item(17, 6).Value = "2026-02-27"
$ws.Cells.Item(17, 7).Value = "https://www.indeed.com/viewjob?jk=e21305a8178b20fc"

$ws.Cells.Item(18, 1).Value = "Home Loans Senior Data Analyst"
$ws.Cells.Item(18, 2).Value = "SoFi"
$ws.Cells.Item(18, 3).Value = "Charlotte, NC, US USA"
$ws.Cells.Item(18, 4).Value = 11.1
$ws.Cells.Item(18, 5).Value = "RAG, Redshift, BigQuery, Git, Snowflake, BigQuery, Redshift, R, Scala, A/B Testing"
$ws.Cells.Item(18, 6).Value = "2026-02-27"
$ws.Cells.Item(18, 7).Value = "https://www.indeed.com/viewjob?jk=53d81a40d07ae30f"

$ws.Cells.Item(19, 1).Value = "GEN AI Engineer"
$ws.Cells.Item(19, 2).Value = "Realign"
$ws.Cells.Item(19, 3).Value = "Tampa, FL, US USA"
$ws.Cells.Item(19, 4).Value = 10
$ws.Cells.Item(19, 5).Value = "AI Engineer, LangChain, RAG, Prompt Engineering, Python, R, Java, Scala, Optimization"
$ws.Cells.Item(19, 6).Value = "2026-02-27"
$ws.Cells.Item(19, 7).Value = "https://www.indeed.com/viewjob?jk=08fb5dc63f69fdde"

$ws.Cells.Item(20, 1).Value = "GEN AI Engineer"
$ws.Cells.Item(20, 2).Value = "Realign"
$ws.Cells.Item(20, 3).Value = "Dallas, TX, US USA"
$ws.Cells.Item(20, 4).Value = 10
$ws.Cells.Item(20, 5).Value = "AI Engineer, LangChain, RAG, Prompt Engineering, Python, R, Java, Scala, Optimization"
$ws.Cells.Item(20, 6).Value = "2026-02-27"
$ws.Cells.Item(20, 7).Value = "https://www.indeed.com/viewjob?jk=211a24837c2308d1"

$ws.Cells.Item(21, 1).Value = "CloudOps Specialist"
$ws.Cells.Item(21, 2).Value = "TechnoMile"
$ws.Cells.Item(21, 3).Value = "McLean, VA, US USA"
$ws.Cells.Item(21, 4).Value = 10
$ws.Cells.Item(21, 5).Value = "S3, EC2, Kubernetes, CI/CD, Terraform, Python, R, Scala, Optimization"
$ws.Cells.Item(21, 6).Value = "2026-02-27"
$ws.Cells.Item(21, 7).Value = "https://www.indeed.com/viewjob?jk=0602f10e5c420065"

$ws.Cells.Item(22, 1).Value = "Senior Linux Firmware Engineer"
$ws.Cells.Item(22, 2).Value = "Satcon"
$ws.Cells.Item(22, 3).Value = "Irving, TX, US USA"
$ws.Cells.Item(22, 4).Value = 10
$ws.Cells.Item(22, 5).Value = "RAG, Jenkins, Git, MySQL, Python, SQL, R, Java, Optimization"
$ws.Cells.Item(22, 6).Value = "2026-02-27"
$ws.Cells.Item(22, 7).Value = "https://www.indeed.com/viewjob?jk=2ad4f125c916fce4"

$ws.Cells.Item(23, 1).Value = "Applied Machine Learning Scientist I (US)"
$ws.Cells.Item(23, 2).Value = "TD Bank"
$ws.Cells.Item(23, 3).Value = "New York, NY, US USA"
$ws.Cells.Item(23, 4).Value = 10
$ws.Cells.Item(23, 5).Value = "Generative AI, LangChain, RAG, PyTorch, PySpark, Hadoop, Python, R, Scala"
$ws.Cells.Item(23, 6).Value = "2026-02-27"
$ws.Cells.Item(23, 7).Value = "https://www.indeed.com/viewjob?jk=75ccc609b71b23c4"
